$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the task distribution numbers for columns C (OOP) and D (Console MVC apps)
$values = @{
    2  = @(5, 5)
    3  = @(3, 5)
    4  = @(4, 1)
    6  = @(5, 2)
    7  = @(3, 6)
    8  = @(1, 2)
    9  = @(2, 4)
    10 = @(2, 6)
    11 = @(4, 7)
    12 = @(6, 4)
    13 = @(1, 7)
    14 = @(7, 3)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 3).Value = $pair[0]
    $ws.Cells.Item($row, 4).Value = $pair[1]
}

# Update the active selection to C15 (matches diff's selection change from B15 to C15)
$ws.Range("C15").Select()
